$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the selection/view on the existing "m06" sheet: the user
#    had scrolled/selected the whole used range (A1:F57) before
#    moving on to build the new mission sheet.
# ------------------------------------------------------------------
$m06 = $wb.Worksheets.Item("m06")
$m06.Activate() | Out-Null
$m06.Range("A1:F57").Select() | Out-Null

# ------------------------------------------------------------------
# 2. Add the new "m08" sheet after the last existing sheet ("m06").
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "m08"

# Column widths (characters) matching the source sheet's autofit widths.
$newSheet.Columns.Item(1).ColumnWidth = 37.75
$newSheet.Columns.Item(2).ColumnWidth = 23.92

# Header row - reuses the same header labels as the other mission sheets.
$newSheet.Range("A1").Value = "name"
$newSheet.Range("B1").Value = "crc"
$newSheet.Range("C1").Value = "group"
$newSheet.Range("D1").Value = "string id"
$newSheet.Range("E1").Value = "txt eng"
$newSheet.Range("F1").Value = "txt ru"

# Voice-line rows for mission 8 ("Lair Escape Race" - Darcy, female).
$newSheet.Range("A5").Value = "dx_m08_5030_darcy"
$newSheet.Range("A4").Value = "dx_m08_5020_darcy"
$newSheet.Range("A2").Value = "dx_m08_5000_darcy"
$newSheet.Range("A3").Value = "dx_m08_5010_darcy"

$newSheet.Range("B5").Value = "0xa2cc80ca"
$newSheet.Range("B4").Value = "0x82cc84cb"
$newSheet.Range("B2").Value = "0x82ccccc9"
$newSheet.Range("B3").Value = "0xa2ccc8c8"

$newSheet.Range("C2").Value = "female"
$newSheet.Range("C3").Value = "female"
$newSheet.Range("C4").Value = "female"
$newSheet.Range("C5").Value = "female"

# Leftover formatted-but-empty cell at F54 (justify + vertical-center),
# a remnant of formatting carried down from the sheet this was built from.
$newSheet.Range("F54").HorizontalAlignment = -4130
$newSheet.Range("F54").VerticalAlignment = -4108

# Selection on the new sheet: A2:A5 with A2 the active cell.
$newSheet.Range("A2:A5").Select() | Out-Null
